$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "43.435.16",
    # "1.00", trailing zeros, etc.) are not coerced into real numbers,
    # matching the inlineStr/text cells in the source workbook.
    $r.NumberFormat = "@"
    $r.Value = $value
    # Reset style back to Normal so we do not leave a stray cell style
    # (the diff only changes cell text, not formatting).
    $r.Style = "Normal"
}

Set-TextValue "D2" "43.435.16"
Set-TextValue "E2" "  +2.89%  "
Set-TextValue "D3" "2.310.56"
Set-TextValue "E3" "  +1.84%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "311.14"
Set-TextValue "E5" "  +1.49%  "
Set-TextValue "D6" "103.15"
Set-TextValue "E6" "  +5.90%  "
Set-TextValue "E7" "  +1.74%  "
Set-TextValue "E8" "  -0.01%  "
Set-TextValue "E9" "  +8.44%  "
Set-TextValue "D10" "35.80"
Set-TextValue "E10" "  +1.55%  "
Set-TextValue "D11" "0.0815"
Set-TextValue "E11" "  +3.26%  "
Set-TextValue "E12" "  -1.06%  "
Set-TextValue "D13" "7.05"
Set-TextValue "E13" "  +2.69%  "
Set-TextValue "D14" "2.668.16"
Set-TextValue "E14" "  +1.74%  "
Set-TextValue "D15" "15.05"
Set-TextValue "E15" "  +2.47%  "
Set-TextValue "D16" "2.309.16"
Set-TextValue "E16" "  +2.18%  "
Set-TextValue "D17" "0.811"
Set-TextValue "E17" "  +2.60%  "
Set-TextValue "D18" "43.335.15"
Set-TextValue "D19" "12.29"
Set-TextValue "E19" "  +0.23%  "
Set-TextValue "D20" "0.0₃0935"
Set-TextValue "E20" "  +3.41%  "
Set-TextValue "D21" "6.19"
Set-TextValue "E21" "  +3.13%  "
Set-TextValue "D22" "68.18"
Set-TextValue "E22" "  +0.69%  "
Set-TextValue "D23" "241.70"
Set-TextValue "E23" "  +1.94%  "
Set-TextValue "E24" "  +1.61%  "
Set-TextValue "E26" "  +0.03%  "
Set-TextValue "D27" "24.98"
Set-TextValue "E27" "  +6.20%  "
Set-TextValue "E28" "  +8.12%  "
Set-TextValue "D29" "36.90"
Set-TextValue "E29" "  -0.87%  "
Set-TextValue "E30" "  +0.98%  "
Set-TextValue "D31" "171.44"
Set-TextValue "E31" "  +5.15%  "
Set-TextValue "E32" "  +0.66%  "
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  +0.00%  "
Set-TextValue "E34" "  +6.97%  "
Set-TextValue "D35" "17.79"
Set-TextValue "E35" "  +0.48%  "
Set-TextValue "D36" "0.0743"
Set-TextValue "E36" "  +1.17%  "
Set-TextValue "E37" "  -2.12%  "
Set-TextValue "B38" "ARBITRUM"
Set-TextValue "C38" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D38" "1.89"
Set-TextValue "E38" "  +3.57%  "
Set-TextValue "B39" "Kaspa"
Set-TextValue "C39" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.106"
Set-TextValue "E39" "  +1.46%  "
Set-TextValue "E40" "  +1.71%  "
Set-TextValue "E41" "  +5.16%  "
Set-TextValue "D42" "2.30"
Set-TextValue "E42" "  -1.30%  "
Set-TextValue "E43" "  +4.49%  "
Set-TextValue "D44" "1.975.03"
Set-TextValue "E44" "  +1.04%  "
Set-TextValue "D45" "19.18"
Set-TextValue "E45" "  +1.21%  "
Set-TextValue "E46" "  +3.29%  "
Set-TextValue "D47" "9.99"
Set-TextValue "E47" "  +0.19%  "
Set-TextValue "D48" "55.63"
Set-TextValue "E48" "  +3.54%  "
Set-TextValue "B49" "Stacks"
Set-TextValue "C49" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D49" "1.60"
Set-TextValue "E49" "  +8.59%  "
Set-TextValue "B50" "HuobiToken"
Set-TextValue "C50" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D50" "2.92"
Set-TextValue "E50" "  +1.54%  "
Set-TextValue "D51" "2.534.84"
Set-TextValue "E51" "  +1.65%  "
